$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1915343915343915
$ws.Range("C2").Value = 0.5661375661375662
$ws.Range("J2").Value = 0.02328042328042328
$ws.Range("P2").Value = 0.1322751322751323
$ws.Range("S2").Value = 0.08677248677248678
$ws.Range("B3").Value = 0.0106951871657754
$ws.Range("C3").Value = 0.03208556149732621
$ws.Range("J3").Value = 0.0267379679144385
$ws.Range("P3").Value = 0.7522281639928698
$ws.Range("S3").Value = 0.17825311942959
$ws.Range("J4").Value = 0.03333333333333333
$ws.Range("P4").Value = 0.72
$ws.Range("S4").Value = 0.2466666666666667
$ws.Range("B6").Value = 0.06211180124223602
$ws.Range("D6").Value = 0.01397515527950311
$ws.Range("F6").Value = 0.08074534161490683
$ws.Range("J6").Value = 0.234472049689441
$ws.Range("O6").Value = 0.01863354037267081
$ws.Range("Q6").Value = 0.1630434782608696
$ws.Range("R6").Value = 0.06366459627329192
$ws.Range("S6").Value = 0.3633540372670808
$ws.Range("B7").Value = 0.107890499194847
$ws.Range("D7").Value = 0.02415458937198068
$ws.Range("E7").Value = 0.001610305958132045
$ws.Range("F7").Value = 0.05475040257648953
$ws.Range("J7").Value = 0.1384863123993559
$ws.Range("O7").Value = 0.02254428341384863
$ws.Range("Q7").Value = 0.1932367149758454
$ws.Range("R7").Value = 0.07246376811594203
$ws.Range("S7").Value = 0.3848631239935588
$ws.Range("B8").Value = 0.09615384615384616
$ws.Range("D8").Value = 0.01257396449704142
$ws.Range("F8").Value = 0.05917159763313609
$ws.Range("J8").Value = 0.1116863905325444
$ws.Range("O8").Value = 0.02071005917159763
$ws.Range("Q8").Value = 0.1767751479289941
$ws.Range("R8").Value = 0.1079881656804734
$ws.Range("S8").Value = 0.4149408284023668
$ws.Range("B9").Value = 0.09355828220858896
$ws.Range("D9").Value = 0.02760736196319018
$ws.Range("E9").Value = 0.001533742331288344
$ws.Range("F9").Value = 0.05521472392638037
$ws.Range("J9").Value = 0.1441717791411043
$ws.Range("O9").Value = 0.01533742331288344
$ws.Range("Q9").Value = 0.1733128834355828
$ws.Range("R9").Value = 0.07515337423312883
$ws.Range("S9").Value = 0.4141104294478528
$ws.Range("B10").Value = 0.1082486995293535
$ws.Range("D10").Value = 0.02254149120634134
$ws.Range("E10").Value = 0.0009908347783007183
$ws.Range("F10").Value = 0.05870696061431756
$ws.Range("J10").Value = 0.1387168689621006
$ws.Range("O10").Value = 0.01709189992568739
$ws.Range("Q10").Value = 0.2142680208075304
$ws.Range("R10").Value = 0.08100074312608373
$ws.Range("S10").Value = 0.3584344810502849
$ws.Range("G11").Value = 0.1451965065502183
$ws.Range("J11").Value = 0.07860262008733625
$ws.Range("K11").Value = 0.1910480349344978
$ws.Range("L11").Value = 0.5764192139737991
$ws.Range("S11").Value = 0.008733624454148471
$ws.Range("F12").Value = 0.001811594202898551
$ws.Range("G12").Value = 0.7481884057971014
$ws.Range("J12").Value = 0.1721014492753623
$ws.Range("K12").Value = 0.007246376811594203
$ws.Range("L12").Value = 0.03985507246376811
$ws.Range("S12").Value = 0.03079710144927536
$ws.Range("F15").Value = 0.02489019033674963
$ws.Range("H15").Value = 0.1376281112737921
$ws.Range("I15").Value = 0.06734992679355783
$ws.Range("J15").Value = 0.3748169838945827
$ws.Range("K15").Value = 0.0746705710102489
$ws.Range("M15").Value = 0.01610541727672035
$ws.Range("O15").Value = 0.08345534407027819
$ws.Range("S15").Value = 0.2210834553440703
$ws.Range("F16").Value = 0.0144
$ws.Range("H16").Value = 0.184
$ws.Range("I16").Value = 0.0784
$ws.Range("J16").Value = 0.4336
$ws.Range("K16").Value = 0.0992
$ws.Range("M16").Value = 0.0176
$ws.Range("N16").Value = 0.0016
$ws.Range("O16").Value = 0.0608
$ws.Range("S16").Value = 0.1104
$ws.Range("F17").Value = 0.01678321678321678
$ws.Range("H17").Value = 0.1727272727272727
$ws.Range("I17").Value = 0.1
$ws.Range("J17").Value = 0.4034965034965035
$ws.Range("K17").Value = 0.0972027972027972
$ws.Range("M17").Value = 0.01608391608391608
$ws.Range("N17").Value = 0.001398601398601399
$ws.Range("O17").Value = 0.06993006993006994
$ws.Range("S17").Value = 0.1223776223776224
$ws.Range("F18").Value = 0.02283849918433932
$ws.Range("H18").Value = 0.2022838499184339
$ws.Range("I18").Value = 0.07340946166394779
$ws.Range("J18").Value = 0.432300163132137
$ws.Range("K18").Value = 0.09298531810766721
$ws.Range("M18").Value = 0.02283849918433932
$ws.Range("N18").Value = 0.001631321370309951
$ws.Range("O18").Value = 0.04404567699836868
$ws.Range("S18").Value = 0.1076672104404568
$ws.Range("F19").Value = 0.01652465788794216
$ws.Range("H19").Value = 0.1995868835528014
$ws.Range("I19").Value = 0.09450038729666925
$ws.Range("J19").Value = 0.3746449780531887
$ws.Range("K19").Value = 0.1107668474051123
$ws.Range("M19").Value = 0.01962303124193132
$ws.Range("N19").Value = 0.0002581977794990963
$ws.Range("O19").Value = 0.0630002581977795
$ws.Range("S19").Value = 0.1210947585850762
